$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '52.043.74'
$ws.Range("E2").Value = '  -0.57%  '
$ws.Range("D3").Value = '2.919.49'
$ws.Range("E3").Value = '  +0.43%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '356.64'
$ws.Range("E5").Value = '  +1.30%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '109.83'
$ws.Range("E6").Value = '  -1.72%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.568'
$ws.Range("E7").Value = '  +1.99%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.634'
$ws.Range("E9").Value = '  +1.07%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.93'
$ws.Range("E10").Value = '  -2.34%  '
$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0873'
$ws.Range("E11").Value = '  +1.02%  '
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.138'
$ws.Range("E12").Value = '  +1.56%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.51'
$ws.Range("E13").Value = '  -1.65%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.78'
$ws.Range("E14").Value = '  +0.03%  '
$ws.Range("D15").Value = '3.372.30'
$ws.Range("E15").Value = '  +0.28%  '
$ws.Range("D16").Value = '2.906.39'
$ws.Range("E16").Value = '  -0.07%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.988'
$ws.Range("E17").Value = '  -2.02%  '
$ws.Range("D18").Value = '51.996.85'
$ws.Range("E18").Value = '  -0.57%  '
$ws.Range("E19").Value = '  +3.74%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.54'
$ws.Range("E20").Value = '  -1.07%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.91'
$ws.Range("E21").Value = '  -1.87%  '
$ws.Range("D22").Value = '0.0₃0980'
$ws.Range("E22").Value = '  +0.19%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.52'
$ws.Range("E23").Value = '  -0.34%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '268.72'
$ws.Range("E24").Value = '  -0.27%  '
$ws.Range("E25").Value = '  +1.88%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.184'
$ws.Range("E26").Value = '  +9.63%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '26.94'
$ws.Range("E27").Value = '  +0.83%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.66'
$ws.Range("E28").Value = '  +15.86%  '
$ws.Range("E29").Value = '  +0.14%  '
$ws.Range("E30").Value = '  +7.15%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '10.49'
$ws.Range("E31").Value = '  -1.20%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '37.69'
$ws.Range("E32").Value = '  +0.11%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.21'
$ws.Range("E33").Value = '  -1.75%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.22'
$ws.Range("E34").Value = '  -1.21%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '52.23'
$ws.Range("E35").Value = '  -1.98%  '
$ws.Range("E36").Value = '  -1.43%  '
$ws.Range("E37").Value = '  +0.05%  '
$ws.Range("E38").Value = '  -2.89%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.24'
$ws.Range("E39").Value = '  -2.70%  '
$ws.Range("E40").Value = '  -3.27%  '
$ws.Range("E41").Value = '  -4.47%  '
$ws.Range("E42").Value = '  +2.58%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '23.00'
$ws.Range("E43").Value = '  -2.39%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '119.62'
$ws.Range("E44").Value = '  -0.79%  '
$ws.Range("E45").Value = '  -0.88%  '
$ws.Range("B46").Value = 'NEARProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.46'
$ws.Range("E46").Value = '  -2.31%  '
$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.47'
$ws.Range("E47").Value = '  -6.32%  '
$ws.Range("D48").Value = '2.126.29'
$ws.Range("E48").Value = '  -3.16%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.252'
$ws.Range("E49").Value = '  -5.04%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0340'
$ws.Range("E50").Value = '  +1.31%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.923'
$ws.Range("E51").Value = '  -4.53%  '
